$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's results (daily update) as row 73
$row = 73

$cellA = $ws.Cells.Item($row, 1)
$cellA.Value = 46022
$cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 2).Value = 162
$ws.Cells.Item($row, 3).Value = 172
$ws.Cells.Item($row, 4).Value = 162
